# Update countries & provincias Spain
# Refreshes the "Pais" COVID snapshot sheet to the next data pull:
#   - timestamp bumped from 05:09 to 06:26
#   - Kazajistan overtakes Ecuador and Suecia (rows 30-32 reshuffle)
#   - Haiti overtakes Guayana Francesa (rows 90-91 reshuffle)
#   - Groenlandia overtakes Islas Malvinas (rows 210-211, tied totals)
#   - plain data refresh for Belgica, Honduras, Australia, Montenegro,
#     Mongolia and Belice rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 25 de Julio de 2020 a las 06:26"

# --- Ranking reshuffle: Kazajistan now above Ecuador / Suecia --------
$ws.Range("A30").Value = "Kazajistan"
$ws.Range("B30").Value = 80226
$ws.Range("C30").Value = 1740
$ws.Range("D30").Value = 51260
$ws.Range("E30").Value = 28381
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 585

$ws.Range("A31").Value = "Ecuador"
$ws.Range("B31").Value = 79049
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 34544
$ws.Range("E31").Value = 39037
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 5468

$ws.Range("A32").Value = "Suecia"
$ws.Range("B32").Value = 78997
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 5697

# --- Plain data refresh: Belgica (row 37) -----------------------------
$ws.Range("B37").Value = 65199
$ws.Range("C37").Value = 352
$ws.Range("D37").Value = 17394
$ws.Range("E37").Value = 37988
$ws.Range("G37").Value = 5
$ws.Range("H37").Value = 9817

# --- Plain data refresh: Honduras (row 52) ----------------------------
$ws.Range("B52").Value = 37559
$ws.Range("C52").Value = 657
$ws.Range("D52").Value = 4607
$ws.Range("E52").Value = 31891
$ws.Range("G52").Value = 50
$ws.Range("H52").Value = 1061

# --- Plain data refresh: Australia (row 74) ---------------------------
$ws.Range("D74").Value = 8929
$ws.Range("E74").Value = 4874

# --- Ranking reshuffle: Haiti now above Guayana Francesa --------------
$ws.Range("A90").Value = "Haiti"
$ws.Range("B90").Value = 7260
$ws.Range("C90").Value = 63
$ws.Range("D90").Value = 4236
$ws.Range("E90").Value = 2868
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 2
$ws.Range("H90").Value = 156

$ws.Range("A91").Value = "Guayana Francesa"
$ws.Range("B91").Value = 7251
$ws.Range("C91").Value = 0
$ws.Range("D91").Value = 5522
$ws.Range("E91").Value = 1688
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 41

# --- Plain data refresh: Montenegro (row 115) -------------------------
$ws.Range("D115").Value = 643
$ws.Range("E115").Value = 1979

# --- Plain data refresh: Mongolia (row 171) ----------------------------
$ws.Range("D171").Value = 218
$ws.Range("E171").Value = 70

# --- Plain data refresh: Belice (row 194) ------------------------------
$ws.Range("B194").Value = 48
$ws.Range("C194").Value = 1
$ws.Range("D194").Value = 26

# --- Ranking reshuffle: Groenlandia now above Islas Malvinas (tied) ---
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"
